# Scheduled runner update: refresh cached market-price-derived figures
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H-N)
# on a handful of leve rows across each crafting-class sheet. A few rows
# lose their trailing LeveProfit cell entirely when that side of the
# recipe no longer prices out (handled via ClearContents below).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2666
$ws.Range("J19").Value = 1499
$ws.Range("L19").Value = 1499
$ws.Range("N19").Value = -1849

$ws.Range("H41").Value = 302.46667
$ws.Range("J41").Value = 366.66666
$ws.Range("L41").Value = 366.66666
$ws.Range("N41").Value = -1246.66666

$ws.Range("H51").Value = 8342.182000000001
$ws.Range("I51").Value = 11162.2
$ws.Range("K51").Value = 11162.2
$ws.Range("M51").Value = -10678.2

$ws.Range("H55").Value = 562.55
$ws.Range("I55").Value = 406.08334
$ws.Range("J55").Value = 797.25
$ws.Range("K55").Value = 406.08334
$ws.Range("L55").Value = 797.25
$ws.Range("M55").Value = -192.08334
$ws.Range("N55").Value = -1225.25

$ws.Range("H92").Value = 63444.562
$ws.Range("I92").Value = 125876.25
$ws.Range("J92").Value = 1012.875
$ws.Range("K92").Value = 125876.25
$ws.Range("L92").Value = 1012.875
$ws.Range("M92").Value = -124628.25
$ws.Range("N92").Value = -3508.875

$ws.Range("H115").Value = 1249.5
$ws.Range("I115").Value = 1049.2
$ws.Range("J115").Value = 1360.7778
$ws.Range("K115").Value = 3147.6
$ws.Range("L115").Value = 4082.3334
$ws.Range("M115").Value = -1580.6
$ws.Range("N115").Value = -7216.3334

$ws.Range("H138").Value = 3453.6897
$ws.Range("I138").Value = 2529.1738
$ws.Range("J138").Value = 6997.6665
$ws.Range("K138").Value = 7587.5214
$ws.Range("L138").Value = 20992.9995
$ws.Range("M138").Value = -2447.5214
$ws.Range("N138").Value = -31272.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 721.4375
$ws.Range("I2").Value = 496.9091
$ws.Range("J2").Value = 1215.4
$ws.Range("K2").Value = 496.9091
$ws.Range("L2").Value = 1215.4
$ws.Range("M2").Value = -383.9091
$ws.Range("N2").Value = -1441.4

$ws.Range("H32").Value = 19536.285
$ws.Range("I32").Value = 20103.139
$ws.Range("K32").Value = 20103.139
$ws.Range("M32").Value = -19816.139

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H45").Value = 3534.3125
$ws.Range("I45").Value = 2707.7727
$ws.Range("K45").Value = 2707.7727
$ws.Range("M45").Value = -2330.7727

$ws.Range("H63").Value = 2850
$ws.Range("I63").Value = 2850
$ws.Range("K63").Value = 2850
$ws.Range("M63").Value = -2164

$ws.Range("H66").Value = 2850
$ws.Range("I66").Value = 2850
$ws.Range("K66").Value = 14250
$ws.Range("M66").Value = -10818

$ws.Range("H88").Value = 2100.7058
$ws.Range("J88").Value = 2475.5833
$ws.Range("L88").Value = 2475.5833
$ws.Range("N88").Value = -3287.5833

$ws.Range("H91").Value = 2100.7058
$ws.Range("J91").Value = 2475.5833
$ws.Range("L91").Value = 2475.5833
$ws.Range("N91").Value = -5283.5833

$ws.Range("H116").Value = 721.4375
$ws.Range("I116").Value = 496.9091
$ws.Range("J116").Value = 1215.4
$ws.Range("K116").Value = 496.9091
$ws.Range("L116").Value = 1215.4
$ws.Range("M116").Value = 1797.0909
$ws.Range("N116").Value = -5803.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 721.4375
$ws.Range("I3").Value = 496.9091
$ws.Range("J3").Value = 1215.4
$ws.Range("K3").Value = 496.9091
$ws.Range("L3").Value = 1215.4
$ws.Range("M3").Value = -382.9091
$ws.Range("N3").Value = -1443.4

$ws.Range("H86").Value = 3133.1333
$ws.Range("I86").Value = 2755.4
$ws.Range("J86").Value = 3888.6
$ws.Range("K86").Value = 2755.4
$ws.Range("L86").Value = 3888.6
$ws.Range("M86").Value = -1632.4
$ws.Range("N86").Value = -6134.6

$ws.Range("H89").Value = 3133.1333
$ws.Range("I89").Value = 2755.4
$ws.Range("J89").Value = 3888.6
$ws.Range("K89").Value = 13777
$ws.Range("L89").Value = 19443
$ws.Range("M89").Value = -8161
$ws.Range("N89").Value = -30675

$ws.Range("H99").Value = 29690.459
$ws.Range("I99").Value = 39256.52
$ws.Range("J99").Value = 3862.1
$ws.Range("K99").Value = 39256.52
$ws.Range("L99").Value = 3862.1
$ws.Range("M99").Value = -37758.52
$ws.Range("N99").Value = -6858.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 762.93335
$ws.Range("I16").Value = 841.8461
$ws.Range("K16").Value = 841.8461
$ws.Range("M16").Value = -554.8461

$ws.Range("H31").Value = 1461.2941
$ws.Range("I31").Value = 1431.375
$ws.Range("K31").Value = 1431.375
$ws.Range("M31").Value = -1136.375

$ws.Range("H34").Value = 1461.2941
$ws.Range("I34").Value = 1431.375
$ws.Range("K34").Value = 1431.375
$ws.Range("M34").Value = -1229.375

$ws.Range("H99").Value = 3437.182
$ws.Range("I99").Value = 3439.5
$ws.Range("K99").Value = 3439.5
$ws.Range("M99").Value = -1941.5

$ws.Range("H105").Value = 1286.2727
$ws.Range("I105").Value = 1286.2727
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1286.2727
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 460.7273
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 762.93335
$ws.Range("I113").Value = 841.8461
$ws.Range("K113").Value = 841.8461
$ws.Range("M113").Value = 1328.1539

$ws.Range("H126").Value = 3437.182
$ws.Range("I126").Value = 3439.5
$ws.Range("K126").Value = 10318.5
$ws.Range("M126").Value = -7848.5

$ws.Range("H134").Value = 48959.363
$ws.Range("I134").Value = 69033
$ws.Range("J134").Value = 5944.4287
$ws.Range("K134").Value = 207099
$ws.Range("L134").Value = 17833.2861
$ws.Range("M134").Value = -204564
$ws.Range("N134").Value = -22903.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 4999
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 4999
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 14997
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -15563

$ws.Range("H122").Value = 327
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 19499.5
$ws.Range("J10").Value = 19499.5
$ws.Range("L10").Value = 19499.5
$ws.Range("N10").Value = -19837.5

$ws.Range("H20").Value = 339835
$ws.Range("I20").Value = 339835
$ws.Range("K20").Value = 339835
$ws.Range("M20").Value = -339590

$ws.Range("H24").Value = 1000000
$ws.Range("I24").Value = 1000000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1000000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -999827
$ws.Range("N24").ClearContents()

$ws.Range("H101").Value = 25148.5
$ws.Range("J101").Value = 25148.5
$ws.Range("L101").Value = 25148.5
$ws.Range("N101").Value = -31638.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 46327.64
$ws.Range("I22").Value = 74459.87
$ws.Range("K22").Value = 74459.87
$ws.Range("M22").Value = -74164.87

$ws.Range("H27").Value = 46327.64
$ws.Range("I27").Value = 74459.87
$ws.Range("K27").Value = 74459.87
$ws.Range("M27").Value = -74352.87

$ws.Range("H40").Value = 3658.25
$ws.Range("I40").Value = 3658.25
$ws.Range("K40").Value = 3658.25
$ws.Range("M40").Value = -3522.25

$ws.Range("H82").Value = 2028.2916
$ws.Range("I82").Value = 701.2727
$ws.Range("K82").Value = 701.2727
$ws.Range("M82").Value = -340.2727

$ws.Range("H85").Value = 2028.2916
$ws.Range("I85").Value = 701.2727
$ws.Range("K85").Value = 701.2727
$ws.Range("M85").Value = 546.7273

$ws.Range("H100").Value = 2835.0908
$ws.Range("I100").Value = 2242.4285
$ws.Range("J100").Value = 3872.25
$ws.Range("K100").Value = 2242.4285
$ws.Range("L100").Value = 3872.25
$ws.Range("M100").Value = -1701.4285
$ws.Range("N100").Value = -4954.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 8999.5
$ws.Range("I39").Value = 8000
$ws.Range("J39").Value = 9999
$ws.Range("K39").Value = 8000
$ws.Range("L39").Value = 9999
$ws.Range("M39").Value = -7587
$ws.Range("N39").Value = -10825

$ws.Range("H51").Value = 30000
$ws.Range("I51").Value = 30000
$ws.Range("K51").Value = 30000
$ws.Range("M51").Value = -29490

$ws.Range("H52").Value = 21694.5
$ws.Range("I52").Value = 21694.5
$ws.Range("K52").Value = 21694.5
$ws.Range("M52").Value = -21468.5

$ws.Range("H70").Value = 22552.5
$ws.Range("J70").Value = 22552.5
$ws.Range("L70").Value = 22552.5
$ws.Range("N70").Value = -23182.5

$ws.Range("H73").Value = 22552.5
$ws.Range("J73").Value = 22552.5
$ws.Range("L73").Value = 22552.5
$ws.Range("N73").Value = -24736.5

$ws.Range("H135").Value = 68738.336
$ws.Range("J135").Value = 68738.336
$ws.Range("L135").Value = 68738.336
$ws.Range("N135").Value = -78878.336
